# ARPA logistics metrics workbook maintenance edit:
#  - row 13 gets overwritten with a refreshed copy of row 2's sample record
#    (same values, except the start/end daytime are shifted 5 days later)
#  - the two stale trailing sample rows (14 and 15) are removed
#  - the view/selection + a handful of column widths are tidied up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refresh row 13 from row 2, shifting the two date/time columns (F, G) by +5 days ---
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($col in $cols) {
    $src = $ws.Range($col + "2")
    $dst = $ws.Range($col + "13")
    $v = $src.Value2()
    if ($col -eq "F" -or $col -eq "G") {
        $dst.Value = $v + 5
    } else {
        $dst.Value = $v
    }
}

# --- drop the two trailing rows that are no longer needed ---
$ws.Range("A14:A15").EntireRow.Delete()

# --- tidy up a handful of column widths (O:U) ---
$ws.Range("O1").EntireColumn.ColumnWidth = 7.39
$ws.Range("P1").EntireColumn.ColumnWidth = 6.94
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.05
$ws.Range("R1").EntireColumn.ColumnWidth = 8.39
$ws.Range("S1").EntireColumn.ColumnWidth = 11.05
$ws.Range("T1").EntireColumn.ColumnWidth = 7.61
$ws.Range("U1").EntireColumn.ColumnWidth = 13.5

# --- reset the scroll position and selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("K25").Select()
